$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with the two new columns (Año / Mes) ---
$tbl = $ws.ListObjects.Item("Tabla1")
$tbl.ListColumns.Add() | Out-Null
$tbl.ListColumns.Add() | Out-Null
$ws.Range("I2").Value = "Año"
$ws.Range("J2").Value = "Mes"

# --- Fill the new columns with the year/month values for every data row ---
for ($r = 3; $r -le 24; $r++) {
    $ws.Cells.Item($r, 9).Value = 2025
    $ws.Cells.Item($r, 10).Value = 11
}

# --- Clear formatting on the (now unused) spacer row and the header row ---
$ws.Range("1:2").ClearFormats()

# --- Update the active selection to match the saved view state ---
$ws.Range("C8").Select()
